{"js": "// Update the 20x5 table of arithmetic expressions with new values, in\n// row-major (reading) order, matching the target edit.\nconst NEW_VALUES = [\n  \"15+14=\",\n  \"6+46=\",\n  \"78+4=\",\n  \"84-78=\",\n  \"8+44=\",\n  \"88-63=\",\n  \"6+13=\",\n  \"29+61=\",\n  \"80+8=\",\n  \"8+26=\",\n  \"95-27=\",\n  \"70-34=\",\n  \"70-49=\",\n  \"43+10=\",\n  \"72-17=\",\n  \"68-7=\",\n  \"14+2=\",\n  \"22+68=\",\n  \"49-1=\",\n  \"71-7=\",\n  \"27-8=\",\n  \"51+46=\",\n  \"82-46=\",\n  \"87-82=\",\n  \"94-85=\",\n  \"34+13=\",\n  \"58+0=\",\n  \"23+24=\",\n  \"24+43=\",\n  \"99-46=\",\n  \"23-1=\",\n  \"37+56=\",\n  \"4+42=\",\n  \"53+5=\",\n  \"33-31=\",\n  \"47-25=\",\n  \"26+58=\",\n  \"14+54=\",\n  \"28+11=\",\n  \"64-15=\",\n  \"0+35=\",\n  \"86+2=\",\n  \"52-34=\",\n  \"20+59=\",\n  \"11+30=\",\n  \"13+57=\",\n  \"13+12=\",\n  \"57+38=\",\n  \"88-76=\",\n  \"20+4=\",\n  \"43+13=\",\n  \"20+73=\",\n  \"41+48=\",\n  \"17-9=\",\n  \"73-27=\",\n  \"22+71=\",\n  \"79+13=\",\n  \"11+33=\",\n  \"68-36=\",\n  \"80-49=\",\n  \"26+41=\",\n  \"71-23=\",\n  \"8+68=\",\n  \"0+19=\",\n  \"34+20=\",\n  \"93-85=\",\n  \"19+51=\",\n  \"18+80=\",\n  \"37+28=\",\n  \"16+83=\",\n  \"1+32=\",\n  \"98-42=\",\n  \"50+44=\",\n  \"56+36=\",\n  \"53+35=\",\n  \"16+14=\",\n  \"57-11=\",\n  \"84-25=\",\n  \"18-5=\",\n  \"31+50=\",\n  \"75+15=\",\n  \"38-10=\",\n  \"90-58=\",\n  \"89-65=\",\n  \"53+33=\",\n  \"90-81=\",\n  \"20+54=\",\n  \"96-66=\",\n  \"87-42=\",\n  \"40-34=\",\n  \"8+58=\",\n  \"18+65=\",\n  \"9+90=\",\n  \"44-5=\",\n  \"98-26=\",\n  \"17+45=\",\n  \"57-21=\",\n  \"1+75=\",\n  \"28+6=\",\n  \"84+9=\"\n];\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body.\");\n}\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nlet idx = 0;\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  await context.sync();\n  for (const cell of cells.items) {\n    if (idx >= NEW_VALUES.length) {\n      break;\n    }\n    // Setting `.value` rewrites the cell's Range.Text in place, which keeps\n    // the existing run formatting (font/size) and only changes the text.\n    cell.value = NEW_VALUES[idx];\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 20x5 table of arithmetic expressions with new values, in\n# row-major (reading) order, matching the target edit.\n$newValues = @(\n    '15+14=',\n    '6+46=',\n    '78+4=',\n    '84-78=',\n    '8+44=',\n    '88-63=',\n    '6+13=',\n    '29+61=',\n    '80+8=',\n    '8+26=',\n    '95-27=',\n    '70-34=',\n    '70-49=',\n    '43+10=',\n    '72-17=',\n    '68-7=',\n    '14+2=',\n    '22+68=',\n    '49-1=',\n    '71-7=',\n    '27-8=',\n    '51+46=',\n    '82-46=',\n    '87-82=',\n    '94-85=',\n    '34+13=',\n    '58+0=',\n    '23+24=',\n    '24+43=',\n    '99-46=',\n    '23-1=',\n    '37+56=',\n    '4+42=',\n    '53+5=',\n    '33-31=',\n    '47-25=',\n    '26+58=',\n    '14+54=',\n    '28+11=',\n    '64-15=',\n    '0+35=',\n    '86+2=',\n    '52-34=',\n    '20+59=',\n    '11+30=',\n    '13+57=',\n    '13+12=',\n    '57+38=',\n    '88-76=',\n    '20+4=',\n    '43+13=',\n    '20+73=',\n    '41+48=',\n    '17-9=',\n    '73-27=',\n    '22+71=',\n    '79+13=',\n    '11+33=',\n    '68-36=',\n    '80-49=',\n    '26+41=',\n    '71-23=',\n    '8+68=',\n    '0+19=',\n    '34+20=',\n    '93-85=',\n    '19+51=',\n    '18+80=',\n    '37+28=',\n    '16+83=',\n    '1+32=',\n    '98-42=',\n    '50+44=',\n    '56+36=',\n    '53+35=',\n    '16+14=',\n    '57-11=',\n    '84-25=',\n    '18-5=',\n    '31+50=',\n    '75+15=',\n    '38-10=',\n    '90-58=',\n    '89-65=',\n    '53+33=',\n    '90-81=',\n    '20+54=',\n    '96-66=',\n    '87-42=',\n    '40-34=',\n    '8+58=',\n    '18+65=',\n    '9+90=',\n    '44-5=',\n    '98-26=',\n    '17+45=',\n    '57-21=',\n    '1+75=',\n    '28+6=',\n    '84+9='\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($idx -ge $newValues.Count) {\n            break\n        }\n        $cell = $t.Cell($r, $c)\n        # Assigning Range.Text rewrites the cell's text in place, keeping the\n        # existing run formatting (font/size) and only changing the text.\n        $cell.Range.Text = $newValues[$idx]\n        $idx++\n    }\n}\n"}
